# Insert a new data row for "Moormerland" as the new row 2 of the
# "BLP-URLs" sheet, pushing all existing data rows down by one, and
# refresh the _FilterDatabase defined name / used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the Moormerland record.
$ws.Range("A2").Value = 457014
$ws.Range("B2").Value = "Moormerland"
$ws.Range("C2").Value = 53.314314
$ws.Range("D2").Value = 7.485564
$ws.Range("E2").Value = "http://www.mmld.de/download"
$ws.Range("F2").Value = "http://lkleer.maps.arcgis.com/home/webmap/viewer.html?webmap=e4311f176259429d970921af4cf49ab2"

# The sheet's _xlnm._FilterDatabase defined name encodes the filtered
# range; bump its last row by one to track the newly inserted row.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -eq "BLP-URLs!_FilterDatabase") {
        $n.RefersTo = "='BLP-URLs'!`$A`$1:`$G`$407"
    }
}
